$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "62.926.71"
Set-TextValue $ws.Range("E2") "  +4.54%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.357.64"
Set-TextValue $ws.Range("E3") "  +4.53%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.23%  "

# Row 5
Set-TextValue $ws.Range("D5") "560.16"
Set-TextValue $ws.Range("E5") "  +3.87%  "

# Row 6
Set-TextValue $ws.Range("D6") "153.03"
Set-TextValue $ws.Range("E6") "  +5.26%  "

# Row 7
Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  -0.03%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.531"
Set-TextValue $ws.Range("E8") "  -0.07%  "

# Row 9
Set-TextValue $ws.Range("D9") "7.51"
Set-TextValue $ws.Range("E9") "  +1.91%  "

# Row 10
Set-TextValue $ws.Range("E10") "  +4.19%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.437"
Set-TextValue $ws.Range("E11") "  +0.49%  "

# Row 12
Set-TextValue $ws.Range("D12") "3.945.90"
Set-TextValue $ws.Range("E12") "  +4.74%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.138"
Set-TextValue $ws.Range("E13") "  +0.26%  "

# Row 14
Set-TextValue $ws.Range("D14") "27.10"
Set-TextValue $ws.Range("E14") "  +3.71%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.0000181"
Set-TextValue $ws.Range("E15") "  +3.61%  "

# Row 16
Set-TextValue $ws.Range("D16") "63.100.54"
Set-TextValue $ws.Range("E16") "  +4.68%  "

# Row 17
Set-TextValue $ws.Range("D17") "3.329.59"
Set-TextValue $ws.Range("E17") "  +3.49%  "

# Row 18
Set-TextValue $ws.Range("D18") "6.44"
Set-TextValue $ws.Range("E18") "  +3.02%  "

# Row 19
Set-TextValue $ws.Range("D19") "13.80"
Set-TextValue $ws.Range("E19") "  +4.60%  "

# Row 20
Set-TextValue $ws.Range("D20") "8.40"
Set-TextValue $ws.Range("E20") "  +0.29%  "

# Row 21
Set-TextValue $ws.Range("D21") "387.74"
Set-TextValue $ws.Range("E21") "  +1.19%  "

# Row 22
Set-TextValue $ws.Range("B22") "Dai"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D22") "1.00"
Set-TextValue $ws.Range("E22") "  -0.27%  "

# Row 23
Set-TextValue $ws.Range("B23") "Polygon"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D23") "0.540"
Set-TextValue $ws.Range("E23") "  +1.78%  "

# Row 24
Set-TextValue $ws.Range("D24") "70.45"
Set-TextValue $ws.Range("E24") "  +0.37%  "

# Row 25
Set-TextValue $ws.Range("D25") "0.181"
Set-TextValue $ws.Range("E25") "  +5.77%  "

# Row 26
Set-TextValue $ws.Range("D26") "8.84"
Set-TextValue $ws.Range("E26") "  -0.34%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.0₃0962"
Set-TextValue $ws.Range("E27") "  +5.44%  "

# Row 28
Set-TextValue $ws.Range("E28") "  -0.05%  "

# Row 29
Set-TextValue $ws.Range("D29") "6.60"
Set-TextValue $ws.Range("E29") "  +6.35%  "

# Row 30
Set-TextValue $ws.Range("D30") "1.99"
Set-TextValue $ws.Range("E30") "  +4.29%  "

# Row 31
Set-TextValue $ws.Range("D31") "5.67"
Set-TextValue $ws.Range("E31") "  +3.78%  "

# Row 32
Set-TextValue $ws.Range("D32") "23.04"
Set-TextValue $ws.Range("E32") "  +2.45%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.30"
Set-TextValue $ws.Range("E33") "  +6.09%  "

# Row 34
Set-TextValue $ws.Range("D34") "6.73"
Set-TextValue $ws.Range("E34") "  +1.67%  "

# Row 35
Set-TextValue $ws.Range("B35") "Monero"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D35") "160.63"
Set-TextValue $ws.Range("E35") "  +2.37%  "

# Row 36
Set-TextValue $ws.Range("B36") "ImmutableX"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D36") "1.48"
Set-TextValue $ws.Range("E36") "  +8.61%  "

# Row 37
Set-TextValue $ws.Range("D37") "1.89"
Set-TextValue $ws.Range("E37") "  +11.95%  "

# Row 38
Set-TextValue $ws.Range("D38") "27.06"
Set-TextValue $ws.Range("E38") "  +4.72%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.0745"
Set-TextValue $ws.Range("E39") "  +5.03%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.831.44"
Set-TextValue $ws.Range("E40") "  +1.41%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.0313"
Set-TextValue $ws.Range("E41") "  +9.35%  "

# Row 42
Set-TextValue $ws.Range("D42") "4.30"
Set-TextValue $ws.Range("E42") "  +0.57%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.749"
Set-TextValue $ws.Range("E43") "  +3.69%  "

# Row 44
Set-TextValue $ws.Range("D44") "40.75"
Set-TextValue $ws.Range("E44") "  +1.99%  "

# Row 45
Set-TextValue $ws.Range("D45") "1.04"
Set-TextValue $ws.Range("E45") "  +3.73%  "

# Row 46
Set-TextValue $ws.Range("B46") "RenzoRestakedETH"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue $ws.Range("D46") "3.409.24"
Set-TextValue $ws.Range("E46") "  +4.75%  "

# Row 47
Set-TextValue $ws.Range("B47") "InjectiveProtocol"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "22.16"
Set-TextValue $ws.Range("E47") "  +6.60%  "

# Row 48
Set-TextValue $ws.Range("E48") "  +1.46%  "

# Row 49
Set-TextValue $ws.Range("D49") "6.30"
Set-TextValue $ws.Range("E49") "  +1.73%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.807"
Set-TextValue $ws.Range("E50") "  -0.11%  "

# Row 51
Set-TextValue $ws.Range("D51") "282.86"
Set-TextValue $ws.Range("E51") "  +4.64%  "
